$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.807.08'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.638.10'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.62'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5057'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2583'
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06426'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.40'
$ws.Range("E10").Value = '  +4.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07798'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.271'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.867.34'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.637.55'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5622'
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").Value = '0.0₅7643'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.26'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '25.843.81'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.05'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.375'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.921'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.137'
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.804'
$ws.Range("E25").Value = '  -4.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.04'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("E27").Value = '  -1.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.809'
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.58'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.245'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04950'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.288'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.235'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.570'
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.385'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9041'
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5573'
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.556'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.131.88'
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01567'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9964'
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.480'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8035'
$ws.Range("E43").Value = '  +1.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.95'
$ws.Range("E44").Value = '  +1.45%  '
$ws.Range("D45").Value = '1.777.21'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").Value = '  -6.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.67'
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4277'
$ws.Range("E48").Value = '  -3.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.792'
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05032'
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9988'
$ws.Range("E51").Value = '  -0.54%  '
